$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 305; this pushes the existing row 305
# (and everything below it) down by one, extending the used range
# from A1:R407 to A1:R408.
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with the new weekly price record.
$ws.Cells.Item(305, 1).Value = 9
$ws.Cells.Item(305, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(305, 3).Value = "Metropolitana"
$ws.Cells.Item(305, 4).Value = 44809
$ws.Cells.Item(305, 5).Value = 13
$ws.Cells.Item(305, 6).Value = 100112052
$ws.Cells.Item(305, 7).Value = "Albahaca"
$ws.Cells.Item(305, 8).Value = "Sin especificar"
$ws.Cells.Item(305, 9).Value = "Primera"
$ws.Cells.Item(305, 10).Value = 350
$ws.Cells.Item(305, 11).Value = 4500
$ws.Cells.Item(305, 12).Value = 5000
$ws.Cells.Item(305, 13).Value = 4857
$ws.Cells.Item(305, 14).Value = "$/paquete"
$ws.Cells.Item(305, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(305, 16).Value = 4857
$ws.Cells.Item(305, 17).Value = 1
$ws.Cells.Item(305, 18).Value = "Hortaliza"
